$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"25.8100000000006"
$ws.Range("H2").Value = [double]"1.617811329144126e-16"
$ws.Range("K2").Value = [double]"56.97737745753237"
$ws.Range("L2").Value = "[48.56294177379726, 65.39181314126748]"
$ws.Range("O2").Value = [double]"1.729605565039425"
$ws.Range("P2").Value = "[1.5786581702723472, 1.8805529598065034]"
$ws.Range("S2").Value = [double]"63.14094216745298"
$ws.Range("T2").Value = "[58.02249330656345, 68.25939102834252]"
$ws.Range("W2").Value = [double]"18.70514514514558"
$ws.Range("X2").Value = [double]"18.0850850850855"
$ws.Range("Y2").Value = [double]"19.32520520520566"

$ws.Range("E3").Value = [double]"22.25000000000004"
$ws.Range("H3").Value = [double]"1.617811329144126e-16"
$ws.Range("I3").Value = [double]"0.8838852037910726"
$ws.Range("K3").Value = [double]"56.19481504622792"
$ws.Range("L3").Value = "[44.02295907535071, 68.36667101710512]"
$ws.Range("M3").Value = [double]"6.661338147750939e-16"
$ws.Range("N3").Value = [double]"6.661338147750939e-16"
$ws.Range("O3").Value = [double]"-2.352263568453619"
$ws.Range("P3").Value = "[-2.5661057110403114, -2.1384214258669267]"
$ws.Range("S3").Value = [double]"56.3525245650068"
$ws.Range("T3").Value = "[49.98085550088737, 62.72419362912622]"
$ws.Range("W3").Value = [double]"8.329829829829844"
$ws.Range("X3").Value = [double]"7.572572572572587"
$ws.Range("Y3").Value = [double]"9.087087087087102"
